$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-8 (replaces old rows 2-11)
$data = @(
    @("000124-19-6_group1",  0.098,               0.029,   0.01,    "Unlikely"),
    @("1000401-12-0_group1", 0.121,               -0.137,  -0.106,  "Unlikely"),
    @("000615-74-7_group1",  0.07000000000000001, 0.132,   -0.045,  "Unlikely"),
    @("000098-86-2_group2",  0.105,               0.081,   -0.001,  "Unlikely"),
    @("1000309-13-0_group1", -0.034,              -0.026,  -0.09,   "Unlikely"),
    @("054446-78-5_group1",  0.048,               0.055,   -0.044,  "Unlikely"),
    @("000620-14-4_group1",  0.103,               0.046,   -0.038,  "Unlikely")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}

# Remove old rows 9-11 (they no longer exist in the new data range)
$ws.Range("A9:E11").Delete()
